$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row right after the existing data (row 97 -> new row 98)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($newRow, 2).Value = "temperature"

# Column C holds "25" which Excel would otherwise auto-convert to a number;
# force it to stay text to match the source data (all cells stored as strings).
$valueCell = $ws.Cells.Item($newRow, 3)
$valueCell.NumberFormat = "@"
$valueCell.Value = "25"
$valueCell.NumberFormat = "General"
$valueCell.Style = "Normal"

$ws.Cells.Item($newRow, 4).Value = "N/A"
$ws.Cells.Item($newRow, 5).Value = "N/A"
$ws.Cells.Item($newRow, 6).Value = "N/A"
